$d = $word.ActiveDocument

# Locate the paragraph that ends with the "chat" user story - the new
# bullet items are inserted right after it (and before the trailing
# empty paragraph at the end of the list).
$anchor = $d.Paragraphs(9)

# First new bullet: paragraph formatting (style + numbering) is copied
# automatically from the anchor paragraph by InsertParagraphAfter.
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs(10)
$p1.Range.Text = "Voglio visualizzare i dettagli di un viaggio, i partecipanti, le date e l’itinerario"

# Second new bullet, inserted right after the first one.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(11)
$p2.Range.Text = "Voglio aggiungere una spesa al viaggio"
